$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1792575.1
$ws.Range("I38").Value = 4032404.5
$ws.Range("J38").Value = 711.6
$ws.Range("K38").Value = 12097213.5
$ws.Range("L38").Value = 2134.8
$ws.Range("M38").Value = -12096841.5
$ws.Range("N38").Value = -2878.8
$ws.Range("H64").Value = 114343.336
$ws.Range("J64").Value = 3781.6667
$ws.Range("L64").Value = 3781.6667
$ws.Range("N64").Value = -4277.6667
$ws.Range("H67").Value = 114343.336
$ws.Range("J67").Value = 3781.6667
$ws.Range("L67").Value = 3781.6667
$ws.Range("N67").Value = -5497.6667
$ws.Range("H127").Value = 24392670
$ws.Range("I127").Value = 197
$ws.Range("J127").Value = 25002482
$ws.Range("K127").Value = 591
$ws.Range("L127").Value = 75007446
$ws.Range("M127").Value = 4369
$ws.Range("N127").Value = -75017366
$ws.Range("H129").Value = 811.60974
$ws.Range("J129").Value = 879.05884
$ws.Range("L129").Value = 2637.17652
$ws.Range("N129").Value = -12637.17652

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 30588.883
$ws.Range("I2").Value = 1076.8182
$ws.Range("J2").Value = 84694.336
$ws.Range("K2").Value = 1076.8182
$ws.Range("L2").Value = 84694.336
$ws.Range("M2").Value = -963.8181999999999
$ws.Range("N2").Value = -84920.336
$ws.Range("H32").Value = 28712.857
$ws.Range("I32").Value = 4807.92
$ws.Range("K32").Value = 4807.92
$ws.Range("M32").Value = -4520.92
$ws.Range("H61").Value = 2941.1177
$ws.Range("I61").Value = 1200
$ws.Range("J61").Value = 3314.2144
$ws.Range("K61").Value = 1200
$ws.Range("L61").Value = 3314.2144
$ws.Range("M61").Value = -988
$ws.Range("N61").Value = -3738.2144
$ws.Range("H74").Value = 3400
$ws.Range("I74").Value = 2671.1
$ws.Range("J74").Value = 4062.6365
$ws.Range("K74").Value = 2671.1
$ws.Range("L74").Value = 4062.6365
$ws.Range("M74").Value = -1797.1
$ws.Range("N74").Value = -5810.636500000001
$ws.Range("H77").Value = 3400
$ws.Range("I77").Value = 2671.1
$ws.Range("J77").Value = 4062.6365
$ws.Range("K77").Value = 13355.5
$ws.Range("L77").Value = 20313.1825
$ws.Range("M77").Value = -8987.5
$ws.Range("N77").Value = -29049.1825
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H116").Value = 30588.883
$ws.Range("I116").Value = 1076.8182
$ws.Range("J116").Value = 84694.336
$ws.Range("K116").Value = 1076.8182
$ws.Range("L116").Value = 84694.336
$ws.Range("M116").Value = 1217.1818
$ws.Range("N116").Value = -89282.336
$ws.Range("H136").Value = 2941.1177
$ws.Range("I136").Value = 1200
$ws.Range("J136").Value = 3314.2144
$ws.Range("K136").Value = 3600
$ws.Range("L136").Value = 9942.643199999999
$ws.Range("M136").Value = -1050
$ws.Range("N136").Value = -15042.6432

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 30588.883
$ws.Range("I3").Value = 1076.8182
$ws.Range("J3").Value = 84694.336
$ws.Range("K3").Value = 1076.8182
$ws.Range("L3").Value = 84694.336
$ws.Range("M3").Value = -962.8181999999999
$ws.Range("N3").Value = -84922.336
$ws.Range("H22").Value = 306.66666
$ws.Range("I22").Value = 310
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 310
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = -137
$ws.Range("N22").Value = -646
$ws.Range("H27").Value = 29748.666

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15184.959
$ws.Range("I31").Value = 23179.6
$ws.Range("J31").Value = 2779.4827
$ws.Range("K31").Value = 23179.6
$ws.Range("L31").Value = 2779.4827
$ws.Range("M31").Value = -22884.6
$ws.Range("N31").Value = -3369.4827
$ws.Range("H34").Value = 15184.959
$ws.Range("I34").Value = 23179.6
$ws.Range("J34").Value = 2779.4827
$ws.Range("K34").Value = 23179.6
$ws.Range("L34").Value = 2779.4827
$ws.Range("M34").Value = -22977.6
$ws.Range("N34").Value = -3183.4827
$ws.Range("H105").Value = 1593.1428
$ws.Range("I105").Value = 1606
$ws.Range("J105").Value = 1561
$ws.Range("K105").Value = 1606
$ws.Range("L105").Value = 1561
$ws.Range("M105").Value = 141
$ws.Range("N105").Value = -5055
$ws.Range("H107").Value = 636.0769
$ws.Range("I107").Value = 671.5625
$ws.Range("J107").Value = 579.3
$ws.Range("K107").Value = 671.5625
$ws.Range("L107").Value = 579.3
$ws.Range("M107").Value = 1248.4375
$ws.Range("N107").Value = -4419.3
$ws.Range("H132").Value = 4957.273
$ws.Range("I132").Value = 5652.4
$ws.Range("J132").Value = 4378
$ws.Range("K132").Value = 16957.2
$ws.Range("L132").Value = 13134
$ws.Range("M132").Value = -14427.2
$ws.Range("N132").Value = -18194
$ws.Range("H134").Value = 1401.7142
$ws.Range("I134").Value = 1402
$ws.Range("J134").Value = 1400
$ws.Range("K134").Value = 4206
$ws.Range("L134").Value = 4200
$ws.Range("M134").Value = -1671
$ws.Range("N134").Value = -9270

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2199.9
$ws.Range("J34").Value = 2712.375
$ws.Range("L34").Value = 8137.125
$ws.Range("N34").Value = -8305.125
$ws.Range("H44").Value = 906.25
$ws.Range("I44").Value = 185
$ws.Range("K44").Value = 555
$ws.Range("M44").Value = -157
$ws.Range("H99").Value = 1963.8
$ws.Range("I99").Value = 1868.75
$ws.Range("K99").Value = 5606.25
$ws.Range("M99").Value = -3360.25
$ws.Range("H131").Value = 825.4897999999999
$ws.Range("I131").Value = 432
$ws.Range("J131").Value = 846.64514
$ws.Range("K131").Value = 1296
$ws.Range("L131").Value = 2539.93542
$ws.Range("M131").Value = 3744
$ws.Range("N131").Value = -12619.93542

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2693.5625
$ws.Range("I132").Value = 2033.2222
$ws.Range("J132").Value = 3542.5715
$ws.Range("K132").Value = 6099.6666
$ws.Range("L132").Value = 10627.7145
$ws.Range("M132").Value = -3569.6666
$ws.Range("N132").Value = -15687.7145

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1971.8077
$ws.Range("I7").Value = 1405.1177
$ws.Range("J7").Value = 3042.2222
$ws.Range("K7").Value = 1405.1177
$ws.Range("L7").Value = 3042.2222
$ws.Range("M7").Value = -1293.1177
$ws.Range("N7").Value = -3266.2222
$ws.Range("H46").Value = 2531527.5
$ws.Range("I46").Value = 3200
$ws.Range("J46").Value = 3374303.2
$ws.Range("K46").Value = 3200
$ws.Range("L46").Value = 3374303.2
$ws.Range("M46").Value = -3012
$ws.Range("N46").Value = -3374679.2
$ws.Range("H61").Value = 2722.7273
$ws.Range("I61").Value = 1934
$ws.Range("K61").Value = 1934
$ws.Range("M61").Value = -1732
$ws.Range("H113").Value = 2722.7273
$ws.Range("I113").Value = 1934
$ws.Range("K113").Value = 1934
$ws.Range("M113").Value = 236
$ws.Range("H122").Value = 4393.8335
$ws.Range("I122").Value = 1468
$ws.Range("J122").Value = 7319.6665
$ws.Range("K122").Value = 4404
$ws.Range("L122").Value = 21958.9995
$ws.Range("M122").Value = -1954
$ws.Range("N122").Value = -26858.9995
$ws.Range("H126").Value = 1971.8077
$ws.Range("I126").Value = 1405.1177
$ws.Range("J126").Value = 3042.2222
$ws.Range("K126").Value = 4215.3531
$ws.Range("L126").Value = 9126.6666
$ws.Range("M126").Value = -1745.3531
$ws.Range("N126").Value = -14066.6666

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1227.2084
$ws.Range("J122").Value = 1488.8182
$ws.Range("L122").Value = 4466.4546
$ws.Range("N122").Value = -9366.454600000001
$ws.Range("H126").Value = 2122
$ws.Range("I126").Value = 2152.5
$ws.Range("K126").Value = 6457.5
$ws.Range("M126").Value = -3987.5
